$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ACCU")
$ws.Cells.Item(2, 2).Value = -986487.1946236491
$ws.Cells.Item(2, 4).Value = 650863.4871826086
$ws.Cells.Item(2, 5).Value = -335623.7074410405
$ws.Cells.Item(3, 2).Value = -917981.1394414512
$ws.Cells.Item(3, 4).Value = 584640.2962321137
$ws.Cells.Item(3, 5).Value = -333340.8432093374
$ws.Cells.Item(4, 2).Value = -849475.0842592533
$ws.Cells.Item(4, 4).Value = 519836.4398363263
$ws.Cells.Item(4, 5).Value = -329638.6444229271
$ws.Cells.Item(5, 2).Value = -780969.0290770555
$ws.Cells.Item(5, 4).Value = 456966.9534278983
$ws.Cells.Item(5, 5).Value = -324002.0756491572
$ws.Cells.Item(6, 2).Value = -712462.9738948577
$ws.Cells.Item(6, 4).Value = 396581.832707184
$ws.Cells.Item(6, 5).Value = -315881.1411876737
$ws.Cells.Item(7, 2).Value = -643956.9187126597
$ws.Cells.Item(7, 4).Value = 339224.1823903473
$ws.Cells.Item(7, 5).Value = -304732.7363223124
$ws.Cells.Item(8, 2).Value = -575450.8635304619
$ws.Cells.Item(8, 4).Value = 285387.6709522555
$ws.Cells.Item(8, 5).Value = -290063.1925782064
$ws.Cells.Item(9, 2).Value = -506944.8083482641
$ws.Cells.Item(9, 4).Value = 235479.6799192216
$ws.Cells.Item(9, 5).Value = -271465.1284290425
$ws.Cells.Item(10, 2).Value = -438438.7531660662
$ws.Cells.Item(10, 4).Value = 189794.7847070808
$ws.Cells.Item(10, 5).Value = -248643.9684589854
$ws.Cells.Item(11, 2).Value = -369932.6979838684
$ws.Cells.Item(11, 4).Value = 148500.8332236
$ws.Cells.Item(11, 5).Value = -221431.8647602684
$ws.Cells.Item(12, 2).Value = -301426.6428016705
$ws.Cells.Item(12, 4).Value = 111637.5749940472
$ws.Cells.Item(12, 5).Value = -189789.0678076232
$ws.Cells.Item(13, 2).Value = -232920.5876194726
$ws.Cells.Item(13, 4).Value = 79126.00646850033
$ws.Cells.Item(13, 5).Value = -153794.5811509723
$ws.Cells.Item(14, 2).Value = -164414.5324372748
$ws.Cells.Item(14, 4).Value = 50785.54440606973
$ws.Cells.Item(14, 5).Value = -113628.988031205
$ws.Cells.Item(15, 2).Value = -95908.47725507691
$ws.Cells.Item(15, 4).Value = 26355.80280941712
$ws.Cells.Item(15, 5).Value = -69552.6744456598
$ws.Cells.Item(16, 2).Value = -27402.42207287905
$ws.Cells.Item(16, 4).Value = 5519.973762853325
$ws.Cells.Item(16, 5).Value = -21882.44831002572
$ws.Cells.Item(17, 2).Value = 41103.63310931881
$ws.Cells.Item(17, 4).Value = -12072.61204373687
$ws.Cells.Item(17, 5).Value = 29031.02106558194
$ws.Cells.Item(18, 2).Value = 109609.6882915167
$ws.Cells.Item(18, 4).Value = -26786.44462601388
$ws.Cells.Item(18, 5).Value = 82823.24366550279
$ws.Cells.Item(19, 2).Value = 178115.7434737145
$ws.Cells.Item(19, 4).Value = -38983.29514938649
$ws.Cells.Item(19, 5).Value = 139132.4483243281
$ws.Cells.Item(20, 2).Value = 246621.7986559124
$ws.Cells.Item(20, 4).Value = -49009.54565294661
$ws.Cells.Item(20, 5).Value = 197612.2530029658
$ws.Cells.Item(21, 2).Value = 315127.8538381102
$ws.Cells.Item(21, 4).Value = -57187.29431586179
$ws.Cells.Item(21, 5).Value = 257940.5595222484
$ws.Cells.Item(22, 2).Value = 383633.9090203081
$ws.Cells.Item(22, 4).Value = -63808.83854585071
$ws.Cells.Item(22, 5).Value = 319825.0704744574
$ws.Cells.Item(23, 2).Value = 452139.9642025059
$ws.Cells.Item(23, 4).Value = -69133.97384683558
$ws.Cells.Item(23, 5).Value = 383005.9903556703
$ws.Cells.Item(24, 2).Value = 520646.0193847038
$ws.Cells.Item(24, 4).Value = -73389.49146940511
$ws.Cells.Item(24, 5).Value = 447256.5279152987
$ws.Cells.Item(25, 2).Value = 589152.0745669017
$ws.Cells.Item(25, 4).Value = -76770.27929848063
$ws.Cells.Item(25, 5).Value = 512381.7952684211
$ws.Cells.Item(26, 2).Value = 657658.1297490995
$ws.Cells.Item(26, 4).Value = -79441.49844061692
$ws.Cells.Item(26, 5).Value = 578216.6313084826
$ws.Cells.Item(27, 2).Value = 726164.1849312974
$ws.Cells.Item(27, 4).Value = -81541.39862857957
$ws.Cells.Item(27, 5).Value = 644622.7863027179

$ws = $wb.Worksheets.Item("NZU")
$ws.Cells.Item(2, 2).Value = -1059392.721239566
$ws.Cells.Item(2, 4).Value = 673330.9735207614
$ws.Cells.Item(2, 5).Value = -386061.7477188047
$ws.Cells.Item(3, 2).Value = -1016152.202005298
$ws.Cells.Item(3, 4).Value = 632575.305673066
$ws.Cells.Item(3, 5).Value = -383576.8963322322
$ws.Cells.Item(4, 2).Value = -972911.6827710301
$ws.Cells.Item(4, 4).Value = 592467.8371254243
$ws.Cells.Item(4, 5).Value = -380443.8456456058
$ws.Cells.Item(5, 2).Value = -929671.1635367621
$ws.Cells.Item(5, 4).Value = 553104.8654683344
$ws.Cells.Item(5, 5).Value = -376566.2980684277
$ws.Cells.Item(6, 2).Value = -886430.6443024941
$ws.Cells.Item(6, 4).Value = 514582.375016051
$ws.Cells.Item(6, 5).Value = -371848.2692864431
$ws.Cells.Item(7, 2).Value = -843190.1250682261
$ws.Cells.Item(7, 4).Value = 476993.8134459599
$ws.Cells.Item(7, 5).Value = -366196.3116222663
$ws.Cells.Item(8, 2).Value = -799949.6058339581
$ws.Cells.Item(8, 4).Value = 440427.9988892117
$ws.Cells.Item(8, 5).Value = -359521.6069447464
$ws.Cells.Item(9, 2).Value = -756709.08659969
$ws.Cells.Item(9, 4).Value = 404967.2418624818
$ws.Cells.Item(9, 5).Value = -351741.8447372083
$ws.Cells.Item(10, 2).Value = -713468.5673654221
$ws.Cells.Item(10, 4).Value = 370685.7452515414
$ws.Cells.Item(10, 5).Value = -342782.8221138808
$ws.Cells.Item(11, 2).Value = -670228.048131154
$ws.Cells.Item(11, 4).Value = 337648.3235168048
$ws.Cells.Item(11, 5).Value = -332579.7246143493
$ws.Cells.Item(12, 2).Value = -626987.5288968862
$ws.Cells.Item(12, 4).Value = 305909.4611922852
$ws.Cells.Item(12, 5).Value = -321078.067704601
$ws.Cells.Item(13, 2).Value = -583747.009662618
$ws.Cells.Item(13, 4).Value = 275512.7119094655
$ws.Cells.Item(13, 5).Value = -308234.2977531526
$ws.Cells.Item(14, 2).Value = -540506.49042835
$ws.Cells.Item(14, 4).Value = 246490.4234305817
$ws.Cells.Item(14, 5).Value = -294016.0669977684
$ws.Cells.Item(15, 2).Value = -497265.9711940821
$ws.Cells.Item(15, 4).Value = 218863.7619250052
$ws.Cells.Item(15, 5).Value = -278402.2092690769
$ws.Cells.Item(16, 2).Value = -454025.4519598141
$ws.Cells.Item(16, 4).Value = 192643.0000185419
$ws.Cells.Item(16, 5).Value = -261382.4519412722
$ws.Cells.Item(17, 2).Value = -410784.9327255461
$ws.Cells.Item(17, 4).Value = 167828.0277761018
$ws.Cells.Item(17, 5).Value = -242956.9049494442
$ws.Cells.Item(18, 2).Value = -367544.4134912781
$ws.Cells.Item(18, 4).Value = 144409.0433539189
$ws.Cells.Item(18, 5).Value = -223135.3701373592
$ws.Cells.Item(19, 2).Value = -324303.89425701
$ws.Cells.Item(19, 4).Value = 122367.3800890443
$ws.Cells.Item(19, 5).Value = -201936.5141679657
$ws.Cells.Item(20, 2).Value = -281063.375022742
$ws.Cells.Item(20, 4).Value = 101676.4287541039
$ws.Cells.Item(20, 5).Value = -179386.9462686381
$ws.Cells.Item(21, 2).Value = -237822.855788474
$ws.Cells.Item(21, 4).Value = 82302.61707717703
$ws.Cells.Item(21, 5).Value = -155520.238711297
$ws.Cells.Item(22, 2).Value = -194582.336554206
$ws.Cells.Item(22, 4).Value = 64206.41293508763
$ws.Cells.Item(22, 5).Value = -130375.9236191184
$ws.Cells.Item(23, 2).Value = -151341.817319938
$ws.Cells.Item(23, 4).Value = 47343.32246084004
$ws.Cells.Item(23, 5).Value = -103998.494859098
$ws.Cells.Item(24, 2).Value = -108101.29808567
$ws.Cells.Item(24, 4).Value = 31664.85932142164
$ws.Cells.Item(24, 5).Value = -76436.43876424836
$ws.Cells.Item(25, 2).Value = -64860.77885140201
$ws.Cells.Item(25, 4).Value = 17119.46635344421
$ws.Cells.Item(25, 5).Value = -47741.3124979578
$ws.Cells.Item(26, 2).Value = -21620.259617134
$ws.Cells.Item(26, 4).Value = 3653.3753934573
$ws.Cells.Item(26, 5).Value = -17966.8842236767
$ws.Cells.Item(27, 2).Value = 21620.259617134
$ws.Cells.Item(27, 4).Value = -8788.604627934286
$ws.Cells.Item(27, 5).Value = 12831.65498919971
$ws.Cells.Item(28, 2).Value = 64860.77885140201
$ws.Cells.Item(28, 4).Value = -20262.37752565523
$ws.Cells.Item(28, 5).Value = 44598.40132574679
$ws.Cells.Item(29, 2).Value = 108101.29808567
$ws.Cells.Item(29, 4).Value = -30823.93063493883
$ws.Cells.Item(29, 5).Value = 77277.36745073117
$ws.Cells.Item(30, 2).Value = 151341.817319938
$ws.Cells.Item(30, 4).Value = -40528.7868685981
$ws.Cells.Item(30, 5).Value = 110813.0304513399
$ws.Cells.Item(31, 2).Value = 194582.336554206
$ws.Cells.Item(31, 4).Value = -49431.53627117059
$ws.Cells.Item(31, 5).Value = 145150.8002830354
$ws.Cells.Item(32, 2).Value = 237822.855788474
$ws.Cells.Item(32, 4).Value = -57585.4436346151
$ws.Cells.Item(32, 5).Value = 180237.4121538589
$ws.Cells.Item(33, 2).Value = 281063.375022742
$ws.Cells.Item(33, 4).Value = -65042.12752490039
$ws.Cells.Item(33, 5).Value = 216021.2474978416
$ws.Cells.Item(34, 2).Value = 324303.89425701
$ws.Cells.Item(34, 4).Value = -71851.3052331088
$ws.Cells.Item(34, 5).Value = 252452.5890239012
$ws.Cells.Item(35, 2).Value = 367544.4134912781
$ws.Cells.Item(35, 4).Value = -78060.59765224076
$ws.Cells.Item(35, 5).Value = 289483.8158390373
$ws.Cells.Item(36, 2).Value = 410784.9327255461
$ws.Cells.Item(36, 4).Value = -83715.3878365335
$ws.Cells.Item(36, 5).Value = 327069.5448890126
$ws.Cells.Item(37, 2).Value = 454025.4519598141
$ws.Cells.Item(37, 4).Value = -88858.72697203014
$ws.Cells.Item(37, 5).Value = 365166.7249877839
$ws.Cells.Item(38, 2).Value = 497265.9711940821
$ws.Cells.Item(38, 4).Value = -93531.28162806781
$ws.Cells.Item(38, 5).Value = 403734.6895660143
$ws.Cells.Item(39, 2).Value = 540506.49042835
$ws.Cells.Item(39, 4).Value = -97771.31642751006
$ws.Cells.Item(39, 5).Value = 442735.17400084
$ws.Cells.Item(40, 2).Value = 583747.009662618
$ws.Cells.Item(40, 4).Value = -101614.7066327444
$ws.Cells.Item(40, 5).Value = 482132.3030298736
$ws.Cells.Item(41, 2).Value = 626987.5288968862
$ws.Cells.Item(41, 4).Value = -105094.9755641006
$ws.Cells.Item(41, 5).Value = 521892.5533327856
$ws.Cells.Item(42, 2).Value = 670228.048131154
$ws.Cells.Item(42, 4).Value = -108243.3522222215
$ws.Cells.Item(42, 5).Value = 561984.6959089325
$ws.Cells.Item(43, 2).Value = 713468.5673654221
$ws.Cells.Item(43, 4).Value = -111088.8449556423
$ws.Cells.Item(43, 5).Value = 602379.7224097798
$ws.Cells.Item(44, 2).Value = 756709.08659969
$ws.Cells.Item(44, 4).Value = -113658.3274835201
$ws.Cells.Item(44, 5).Value = 643050.75911617
$ws.Cells.Item(45, 2).Value = 799949.6058339581
$ws.Cells.Item(45, 4).Value = -115976.6340388738
$ws.Cells.Item(45, 5).Value = 683972.9717950843
$ws.Cells.Item(46, 2).Value = 843190.1250682261
$ws.Cells.Item(46, 4).Value = -118066.6608311174
$ws.Cells.Item(46, 5).Value = 725123.4642371088
$ws.Cells.Item(47, 2).Value = 886430.6443024941
$ws.Cells.Item(47, 4).Value = -119949.4714318866
$ws.Cells.Item(47, 5).Value = 766481.1728706075
$ws.Cells.Item(48, 2).Value = 929671.1635367621
$ws.Cells.Item(48, 4).Value = -121644.4040613969
$ws.Cells.Item(48, 5).Value = 808026.7594753653
$ws.Cells.Item(49, 2).Value = 972911.6827710301
$ws.Cells.Item(49, 4).Value = -123169.1790917361
$ws.Cells.Item(49, 5).Value = 849742.503679294
$ws.Cells.Item(50, 2).Value = 1016152.202005298
$ws.Cells.Item(50, 4).Value = -124540.0053880393
$ws.Cells.Item(50, 5).Value = 891612.1966172588
$ws.Cells.Item(51, 2).Value = 1059392.721239566
$ws.Cells.Item(51, 4).Value = -125771.6843787859
$ws.Cells.Item(51, 5).Value = 933621.0368607802
$ws.Cells.Item(52, 2).Value = 1102633.240473834
$ws.Cells.Item(52, 4).Value = -126877.7109838903
$ws.Cells.Item(52, 5).Value = 975755.5294899438
$ws.Cells.Item(53, 2).Value = 1145873.759708102
$ws.Cells.Item(53, 4).Value = -127870.3707355338
$ws.Cells.Item(53, 5).Value = 1018003.388972569
$ws.Cells.Item(54, 2).Value = 1189114.27894237
$ws.Cells.Item(54, 4).Value = -128760.8326041029
$ws.Cells.Item(54, 5).Value = 1060353.446338267
$ws.Cells.Item(55, 2).Value = 1232354.798176638
$ws.Cells.Item(55, 4).Value = -129559.2371925308
$ws.Cells.Item(55, 5).Value = 1102795.560984107
$ws.Cells.Item(56, 2).Value = 1275595.317410906
$ws.Cells.Item(56, 4).Value = -130274.7800892571
$ws.Cells.Item(56, 5).Value = 1145320.537321649
$ws.Cells.Item(57, 2).Value = 1318835.836645174
$ws.Cells.Item(57, 4).Value = -130915.7902755432
$ws.Cells.Item(57, 5).Value = 1187920.046369631
$ws.Cells.Item(58, 2).Value = 1362076.355879442
$ws.Cells.Item(58, 4).Value = -131489.8035693134
$ws.Cells.Item(58, 5).Value = 1230586.552310129
$ws.Cells.Item(59, 2).Value = 1405316.87511371
$ws.Cells.Item(59, 4).Value = -132003.6311574754
$ws.Cells.Item(59, 5).Value = 1273313.243956235
$ws.Cells.Item(60, 2).Value = 1448557.394347978
$ws.Cells.Item(60, 4).Value = -132463.4233239187
$ws.Cells.Item(60, 5).Value = 1316093.971024059
$ws.Cells.Item(61, 2).Value = 1491797.913582246
$ws.Cells.Item(61, 4).Value = -132874.7285231359
$ws.Cells.Item(61, 5).Value = 1358923.18505911
$ws.Cells.Item(62, 2).Value = 1535038.432816514
$ws.Cells.Item(62, 4).Value = -133242.5479814719
$ws.Cells.Item(62, 5).Value = 1401795.884835042

$ws = $wb.Worksheets.Item("EUA")
$ws.Cells.Item(2, 2).Value = -1354932.909375169
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = -1354932.909375169
$ws.Cells.Item(3, 2).Value = -1321543.729745422
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = -1321543.729745422
$ws.Cells.Item(4, 2).Value = -1288154.550115674
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = -1288154.550115674
$ws.Cells.Item(5, 2).Value = -1254765.370485926
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = -1254765.370485926
$ws.Cells.Item(6, 2).Value = -1221376.190856178
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = -1221376.190856178
$ws.Cells.Item(7, 2).Value = -1187987.01122643
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = -1187987.01122643
$ws.Cells.Item(8, 2).Value = -1154597.831596682
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = -1154597.831596682
$ws.Cells.Item(9, 2).Value = -1121208.651966934
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = -1121208.651966934
$ws.Cells.Item(10, 2).Value = -1087819.472337186
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = -1087819.472337186
$ws.Cells.Item(11, 2).Value = -1054430.292707438
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = -1054430.292707438
$ws.Cells.Item(12, 2).Value = -1021041.11307769
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = -1021041.11307769
$ws.Cells.Item(13, 2).Value = -987651.9334479426
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = -987651.9334479426
$ws.Cells.Item(14, 2).Value = -954262.7538181947
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = -954262.7538181947
$ws.Cells.Item(15, 2).Value = -920873.5741884467
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = -920873.5741884467
$ws.Cells.Item(16, 2).Value = -887484.3945586988
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = -887484.3945586988
$ws.Cells.Item(17, 2).Value = -854095.214928951
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = -854095.214928951
$ws.Cells.Item(18, 2).Value = -820706.0352992031
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = -820706.0352992031
$ws.Cells.Item(19, 2).Value = -787316.8556694551
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = -787316.8556694551
$ws.Cells.Item(20, 2).Value = -753927.6760397073
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = -753927.6760397073
$ws.Cells.Item(21, 2).Value = -720538.4964099595
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = -720538.4964099595
$ws.Cells.Item(22, 2).Value = -687149.3167802115
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = -687149.3167802115
$ws.Cells.Item(23, 2).Value = -653760.1371504636
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = -653760.1371504636
$ws.Cells.Item(24, 2).Value = -620370.9575207158
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = -620370.9575207158
$ws.Cells.Item(25, 2).Value = -586981.7778909679
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = -586981.7778909679
$ws.Cells.Item(26, 2).Value = -553592.5982612199
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = -553592.5982612199
$ws.Cells.Item(27, 2).Value = -520203.4186314721
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = -520203.4186314721
$ws.Cells.Item(28, 2).Value = -486814.2390017241
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = -486814.2390017241
$ws.Cells.Item(29, 2).Value = -453425.0593719763
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = -453425.0593719763
$ws.Cells.Item(30, 2).Value = -420035.8797422284
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = -420035.8797422284
$ws.Cells.Item(31, 2).Value = -386646.7001124805
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = -386646.7001124805
$ws.Cells.Item(32, 2).Value = -353257.5204827326
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = -353257.5204827326
$ws.Cells.Item(33, 2).Value = -319868.3408529848
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = -319868.3408529848
$ws.Cells.Item(34, 2).Value = -286479.1612232369
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = -286479.1612232369
$ws.Cells.Item(35, 2).Value = -253089.9815934889
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = -253089.9815934889
$ws.Cells.Item(36, 2).Value = -219700.8019637411
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = -219700.8019637411
$ws.Cells.Item(37, 2).Value = -186311.6223339932
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = -186311.6223339932
$ws.Cells.Item(38, 2).Value = -152922.4427042453
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = -152922.4427042453
$ws.Cells.Item(39, 2).Value = -119533.2630744974
$ws.Cells.Item(39, 4).Value = 0
$ws.Cells.Item(39, 5).Value = -119533.2630744974
$ws.Cells.Item(40, 2).Value = -86144.0834447495
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = -86144.0834447495
$ws.Cells.Item(41, 2).Value = -52754.90381500161
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = -52754.90381500161
$ws.Cells.Item(42, 2).Value = -19365.72418525372
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = -19365.72418525372
$ws.Cells.Item(43, 2).Value = 14023.45544449417
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 14023.45544449417
$ws.Cells.Item(44, 2).Value = 47412.63507424206
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 47412.63507424206
$ws.Cells.Item(45, 2).Value = 80801.81470398995
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 80801.81470398995
$ws.Cells.Item(46, 2).Value = 114190.9943337378
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 114190.9943337378
$ws.Cells.Item(47, 2).Value = 147580.1739634857
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 147580.1739634857
$ws.Cells.Item(48, 2).Value = 180969.3535932336
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 180969.3535932336
$ws.Cells.Item(49, 2).Value = 214358.5332229815
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 214358.5332229815
$ws.Cells.Item(50, 2).Value = 247747.7128527294
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 247747.7128527294
$ws.Cells.Item(51, 2).Value = 281136.8924824773
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 281136.8924824773
$ws.Cells.Item(52, 2).Value = 314526.0721122252
$ws.Cells.Item(52, 4).Value = 0
$ws.Cells.Item(52, 5).Value = 314526.0721122252
$ws.Cells.Item(53, 2).Value = 347915.2517419731
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 347915.2517419731
$ws.Cells.Item(54, 2).Value = 381304.431371721
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 5).Value = 381304.431371721
$ws.Cells.Item(55, 2).Value = 414693.6110014688
$ws.Cells.Item(55, 4).Value = 0
$ws.Cells.Item(55, 5).Value = 414693.6110014688
$ws.Cells.Item(56, 2).Value = 448082.7906312168
$ws.Cells.Item(56, 4).Value = 0
$ws.Cells.Item(56, 5).Value = 448082.7906312168
$ws.Cells.Item(57, 2).Value = 481471.9702609646
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 481471.9702609646
$ws.Cells.Item(58, 2).Value = 514861.1498907126
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 514861.1498907126
$ws.Cells.Item(59, 2).Value = 548250.3295204605
$ws.Cells.Item(59, 4).Value = 0
$ws.Cells.Item(59, 5).Value = 548250.3295204605
$ws.Cells.Item(60, 2).Value = 581639.5091502083
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 581639.5091502083
$ws.Cells.Item(61, 2).Value = 615028.6887799562
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(61, 5).Value = 615028.6887799562
$ws.Cells.Item(62, 2).Value = 648417.8684097041
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = 648417.8684097041
$ws.Cells.Item(63, 2).Value = 681807.0480394519
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 681807.0480394519
$ws.Cells.Item(64, 2).Value = 715196.2276691998
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 715196.2276691998
$ws.Cells.Item(65, 2).Value = 748585.4072989478
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 748585.4072989478
$ws.Cells.Item(66, 2).Value = 781974.5869286957
$ws.Cells.Item(66, 4).Value = 0
$ws.Cells.Item(66, 5).Value = 781974.5869286957
$ws.Cells.Item(67, 2).Value = 815363.7665584435
$ws.Cells.Item(67, 4).Value = 0
$ws.Cells.Item(67, 5).Value = 815363.7665584435
$ws.Cells.Item(68, 2).Value = 848752.9461881914
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 848752.9461881914
$ws.Cells.Item(69, 2).Value = 882142.1258179394
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 882142.1258179394
$ws.Cells.Item(70, 2).Value = 915531.3054476872
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 915531.3054476872
$ws.Cells.Item(71, 2).Value = 948920.485077435
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 948920.485077435
$ws.Cells.Item(72, 2).Value = 982309.664707183
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(72, 5).Value = 982309.664707183

$ws = $wb.Worksheets.Item("UKA")
$ws.Cells.Item(2, 2).Value = -838957.4106661081
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = -838957.4106661081
$ws.Cells.Item(3, 2).Value = -778251.3751620481
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = -778251.3751620481
$ws.Cells.Item(4, 2).Value = -717545.3396579883
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = -717545.3396579883
$ws.Cells.Item(5, 2).Value = -656839.3041539284
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = -656839.3041539284
$ws.Cells.Item(6, 2).Value = -596133.2686498684
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = -596133.2686498684
$ws.Cells.Item(7, 2).Value = -535427.2331458085
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = -535427.2331458085
$ws.Cells.Item(8, 2).Value = -474721.1976417486
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = -474721.1976417486
$ws.Cells.Item(9, 2).Value = -414015.1621376887
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = -414015.1621376887
$ws.Cells.Item(10, 2).Value = -353309.1266336287
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = -353309.1266336287
$ws.Cells.Item(11, 2).Value = -292603.0911295689
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = -292603.0911295689
$ws.Cells.Item(12, 2).Value = -231897.0556255089
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = -231897.0556255089
$ws.Cells.Item(13, 2).Value = -171191.020121449
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = -171191.020121449
$ws.Cells.Item(14, 2).Value = -110484.9846173891
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = -110484.9846173891
$ws.Cells.Item(15, 2).Value = -49778.94911332915
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = -49778.94911332915
$ws.Cells.Item(16, 2).Value = 10927.08639073077
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 10927.08639073077
$ws.Cells.Item(17, 2).Value = 71633.12189479069
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 71633.12189479069
$ws.Cells.Item(18, 2).Value = 132339.1573988506
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 132339.1573988506
$ws.Cells.Item(19, 2).Value = 193045.1929029105
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 193045.1929029105
$ws.Cells.Item(20, 2).Value = 253751.2284069704
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 253751.2284069704
$ws.Cells.Item(21, 2).Value = 314457.2639110303
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 314457.2639110303
$ws.Cells.Item(22, 2).Value = 375163.2994150903
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 375163.2994150903
$ws.Cells.Item(23, 2).Value = 435869.3349191502
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 435869.3349191502
$ws.Cells.Item(24, 2).Value = 496575.3704232101
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 496575.3704232101
$ws.Cells.Item(25, 2).Value = 557281.4059272701
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 557281.4059272701
$ws.Cells.Item(26, 2).Value = 617987.4414313299
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 617987.4414313299
$ws.Cells.Item(27, 2).Value = 678693.4769353899
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 678693.4769353899
$ws.Cells.Item(28, 2).Value = 739399.5124394498
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 739399.5124394498
$ws.Cells.Item(29, 2).Value = 800105.5479435096
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 800105.5479435096
$ws.Cells.Item(30, 2).Value = 860811.5834475696
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 860811.5834475696
$ws.Cells.Item(31, 2).Value = 921517.6189516296
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 921517.6189516296
$ws.Cells.Item(32, 2).Value = 982223.6544556894
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 982223.6544556894
$ws.Cells.Item(33, 2).Value = 1042929.689959749
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 1042929.689959749
$ws.Cells.Item(34, 2).Value = 1103635.725463809
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 1103635.725463809
$ws.Cells.Item(35, 2).Value = 1164341.760967869
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 1164341.760967869
$ws.Cells.Item(36, 2).Value = 1225047.796471929
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = 1225047.796471929
$ws.Cells.Item(37, 2).Value = 1285753.831975989
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = 1285753.831975989
$ws.Cells.Item(38, 2).Value = 1346459.867480049
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 1346459.867480049
$ws.Cells.Item(39, 2).Value = 1407165.902984109
$ws.Cells.Item(39, 4).Value = 0
$ws.Cells.Item(39, 5).Value = 1407165.902984109
$ws.Cells.Item(40, 2).Value = 1467871.938488169
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 1467871.938488169
$ws.Cells.Item(41, 2).Value = 1528577.973992229
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = 1528577.973992229
$ws.Cells.Item(42, 2).Value = 1589284.009496289
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 1589284.009496289
$ws.Cells.Item(43, 2).Value = 1649990.045000348
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 1649990.045000348
$ws.Cells.Item(44, 2).Value = 1710696.080504408
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 1710696.080504408
$ws.Cells.Item(45, 2).Value = 1771402.116008468
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 1771402.116008468
$ws.Cells.Item(46, 2).Value = 1832108.151512528
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 1832108.151512528
$ws.Cells.Item(47, 2).Value = 1892814.187016588
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 1892814.187016588

$ws = $wb.Worksheets.Item("CCA")
$ws.Cells.Item(2, 2).Value = -637370.3534764132
$ws.Cells.Item(2, 4).Value = 357612.2089626566
$ws.Cells.Item(2, 5).Value = -279758.1445137566
$ws.Cells.Item(3, 2).Value = -569129.4162733712
$ws.Cells.Item(3, 4).Value = 305738.6700329768
$ws.Cells.Item(3, 5).Value = -263390.7462403944
$ws.Cells.Item(4, 2).Value = -500888.4790703291
$ws.Cells.Item(4, 4).Value = 256890.6931762101
$ws.Cells.Item(4, 5).Value = -243997.785894119
$ws.Cells.Item(5, 2).Value = -432647.541867287
$ws.Cells.Item(5, 4).Value = 211215.0589002157
$ws.Cells.Item(5, 5).Value = -221432.4829670713
$ws.Cells.Item(6, 2).Value = -364406.6046642449
$ws.Cells.Item(6, 4).Value = 168799.2460639855
$ws.Cells.Item(6, 5).Value = -195607.3586002594
$ws.Cells.Item(7, 2).Value = -296165.6674612028
$ws.Cells.Item(7, 4).Value = 129672.9613602975
$ws.Cells.Item(7, 5).Value = -166492.7061009053
$ws.Cells.Item(8, 2).Value = -227924.7302581607
$ws.Cells.Item(8, 4).Value = 93812.39172674362
$ws.Cells.Item(8, 5).Value = -134112.3385314171
$ws.Cells.Item(9, 2).Value = -159683.7930551187
$ws.Cells.Item(9, 4).Value = 61146.41002558191
$ws.Cells.Item(9, 5).Value = -98537.38302953678
$ws.Cells.Item(10, 2).Value = -91442.85585207661
$ws.Cells.Item(10, 4).Value = 31563.98122344511
$ws.Cells.Item(10, 5).Value = -59878.8746286315
$ws.Cells.Item(11, 2).Value = -23201.91864903453
$ws.Cells.Item(11, 4).Value = 4922.099216976399
$ws.Cells.Item(11, 5).Value = -18279.81943205814
$ws.Cells.Item(12, 2).Value = 45039.01855400753
$ws.Cells.Item(12, 4).Value = -18946.29552444602
$ws.Cells.Item(12, 5).Value = 26092.72302956151
$ws.Cells.Item(13, 2).Value = 113279.9557570496
$ws.Cells.Item(13, 4).Value = -40224.83374285566
$ws.Cells.Item(13, 5).Value = 73055.12201419394
$ws.Cells.Item(14, 2).Value = 181520.8929600917
$ws.Cells.Item(14, 4).Value = -59106.95678303757
$ws.Cells.Item(14, 5).Value = 122413.9361770541
$ws.Cells.Item(15, 2).Value = 249761.8301631338
$ws.Cells.Item(15, 4).Value = -75790.0481033766
$ws.Cells.Item(15, 5).Value = 173971.7820597572
$ws.Cells.Item(16, 2).Value = 318002.7673661758
$ws.Cells.Item(16, 4).Value = -90470.5318184105
$ws.Cells.Item(16, 5).Value = 227532.2355477653
$ws.Cells.Item(17, 2).Value = 386243.7045692179
$ws.Cells.Item(17, 4).Value = -103339.9394046548
$ws.Cells.Item(17, 5).Value = 282903.7651645631
$ws.Cells.Item(18, 2).Value = 454484.64177226
$ws.Cells.Item(18, 4).Value = -114581.8928315436
$ws.Cells.Item(18, 5).Value = 339902.7489407164
$ws.Cells.Item(19, 2).Value = 522725.578975302
$ws.Cells.Item(19, 4).Value = -124369.9179298345
$ws.Cells.Item(19, 5).Value = 398355.6610454675
$ws.Cells.Item(20, 2).Value = 590966.516178344
$ws.Cells.Item(20, 4).Value = -132865.9825410374
$ws.Cells.Item(20, 5).Value = 458100.5336373067
$ws.Cells.Item(21, 2).Value = 659207.4533813862
$ws.Cells.Item(21, 4).Value = -140219.6465023532
$ws.Cells.Item(21, 5).Value = 518987.806879033
$ws.Cells.Item(22, 2).Value = 727448.3905844283
$ws.Cells.Item(22, 4).Value = -146567.711624401
$ws.Cells.Item(22, 5).Value = 580880.6789600273
$ws.Cells.Item(23, 2).Value = 795689.3277874703
$ws.Cells.Item(23, 4).Value = -152034.2667449792
$ws.Cells.Item(23, 5).Value = 643655.0610424911
$ws.Cells.Item(24, 2).Value = 863930.2649905123
$ws.Cells.Item(24, 4).Value = -156731.0334274979
$ws.Cells.Item(24, 5).Value = 707199.2315630144
$ws.Cells.Item(25, 2).Value = 932171.2021935545
$ws.Cells.Item(25, 4).Value = -160757.9301628697
$ws.Cells.Item(25, 5).Value = 771413.2720306847
$ws.Cells.Item(26, 2).Value = 1000412.139396597
$ws.Cells.Item(26, 4).Value = -164203.7857397831
$ws.Cells.Item(26, 5).Value = 836208.3536568135
$ws.Cells.Item(27, 2).Value = 1068653.076599639
$ws.Cells.Item(27, 4).Value = -167147.1448794644
$ws.Cells.Item(27, 5).Value = 901505.9317201743
$ws.Cells.Item(28, 2).Value = 1136894.013802681
$ws.Cells.Item(28, 4).Value = -169657.1207183572
$ws.Cells.Item(28, 5).Value = 967236.8930843235
$ws.Cells.Item(29, 2).Value = 1205134.951005723
$ws.Cells.Item(29, 4).Value = -171794.258945288
$ws.Cells.Item(29, 5).Value = 1033340.692060435
$ws.Cells.Item(30, 2).Value = 1273375.888208765
$ws.Cells.Item(30, 4).Value = -173611.3872209587
$ws.Cells.Item(30, 5).Value = 1099764.500987806
$ws.Cells.Item(31, 2).Value = 1341616.825411807
$ws.Cells.Item(31, 4).Value = -175154.4309171621
$ws.Cells.Item(31, 5).Value = 1166462.394494645
$ws.Cells.Item(32, 2).Value = 1409857.762614849
$ws.Cells.Item(32, 4).Value = -176463.1822844973
$ws.Cells.Item(32, 5).Value = 1233394.580330352
$ws.Cells.Item(33, 2).Value = 1478098.699817891
$ws.Cells.Item(33, 4).Value = -177572.0150139182
$ws.Cells.Item(33, 5).Value = 1300526.684803973
$ws.Cells.Item(34, 2).Value = 1546339.637020933
$ws.Cells.Item(34, 4).Value = -178510.5399472473
$ws.Cells.Item(34, 5).Value = 1367829.097073686
$ws.Cells.Item(35, 2).Value = 1614580.574223975
$ws.Cells.Item(35, 4).Value = -179304.2005703544
$ws.Cells.Item(35, 5).Value = 1435276.373653621
$ws.Cells.Item(36, 2).Value = 1682821.511427017
$ws.Cells.Item(36, 4).Value = -179974.8090412932
$ws.Cells.Item(36, 5).Value = 1502846.702385724
